# Update row 2 attribution values to reflect new relative-direction actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = -0
$ws.Range("B2").Value = -0.09743256086901915
$ws.Range("C2").Value = -0
$ws.Range("D2").Value = 0.2622399841521429
$ws.Range("E2").Value = 0.00303967780147252
$ws.Range("F2").Value = -0
$ws.Range("G2").Value = 0
$ws.Range("I2").Value = -0
$ws.Range("J2").Value = -0
$ws.Range("K2").Value = -0.0157858522608739
$ws.Range("L2").Value = -0
$ws.Range("M2").Value = 0.2500957338607838
$ws.Range("N2").Value = 0.002467049924387657
$ws.Range("R2").Value = -0
$ws.Range("S2").Value = 0
$ws.Range("T2").Value = -0.1032613445660863
$ws.Range("V2").Value = 0.02328870688007144
$ws.Range("W2").Value = -0.04642507234038794
$ws.Range("AB2").Value = 0
$ws.Range("AC2").Value = -0.05884914754576637
$ws.Range("AD2").Value = 0
$ws.Range("AE2").Value = -0.008743753015760445
$ws.Range("AF2").Value = 0.002596583407561699
$ws.Range("AG2").Value = -0
$ws.Range("AI2").Value = -0
$ws.Range("AJ2").Value = 0
$ws.Range("AK2").Value = -0
$ws.Range("AL2").Value = -0.04332689688553018
$ws.Range("AM2").Value = 0
$ws.Range("AN2").Value = 0.03201138138282502
$ws.Range("AO2").Value = 0.07100016630325545
$ws.Range("AQ2").Value = 0
$ws.Range("AR2").Value = -0
$ws.Range("AT2").Value = 0
$ws.Range("AU2").Value = -0.1663206817060559
$ws.Range("AW2").Value = 0.09800526917072176
$ws.Range("AX2").Value = 0.0003144551746630985
$ws.Range("AY2").Value = -0
$ws.Range("BC2").Value = -0
$ws.Range("BD2").Value = -0.02341166372543108
$ws.Range("BF2").Value = 0.1037817066555361
$ws.Range("BG2").Value = 0.03475721051285342
$ws.Range("BJ2").Value = -0
$ws.Range("BL2").Value = 0
$ws.Range("BM2").Value = 0.03593428771602011
$ws.Range("BO2").Value = -0.0417328291051336
$ws.Range("BP2").Value = -0.09108884403239019
$ws.Range("BU2").Value = 0
$ws.Range("BV2").Value = -0.05189477069328403
$ws.Range("BW2").Value = 0
$ws.Range("BX2").Value = 0.02464001639811764
$ws.Range("BY2").Value = -0.02992743296277599
$ws.Range("BZ2").Value = -0
$ws.Range("CB2").Value = 0
$ws.Range("CD2").Value = -0
$ws.Range("CE2").Value = 0.03395303151392331
$ws.Range("CG2").Value = -0.04222034789380898
$ws.Range("CH2").Value = 0.01981383323074086
$ws.Range("CJ2").Value = -0
$ws.Range("CM2").Value = -0
$ws.Range("CN2").Value = -0.01585140079057015
$ws.Range("CP2").Value = 0.0266578807857839
$ws.Range("CQ2").Value = 0.04226796836385542
$ws.Range("CT2").Value = 0
$ws.Range("CU2").Value = -0
$ws.Range("CV2").Value = -0
$ws.Range("CW2").Value = 0.04930545569512913
$ws.Range("CY2").Value = -0.04394995226766105
$ws.Range("CZ2").Value = 0.01160994819653272
$ws.Range("DD2").Value = -0
$ws.Range("DE2").Value = -0
$ws.Range("DF2").Value = 0.03230739597374906
$ws.Range("DH2").Value = 0.01908407623220925
$ws.Range("DI2").Value = 0.04476623474884861
$ws.Range("DJ2").Value = 0
$ws.Range("DK2").Value = -0
$ws.Range("DL2").Value = -0
$ws.Range("DN2").Value = 0
$ws.Range("DO2").Value = -0.02352762696697234
$ws.Range("DQ2").Value = 0.03397899713890748
$ws.Range("DR2").Value = -0.02038357453964971
$ws.Range("DS2").Value = -0
$ws.Range("DW2").Value = 0
$ws.Range("DX2").Value = -0.05521342815995583
$ws.Range("DY2").Value = -0
$ws.Range("DZ2").Value = -0.01748659847869833
$ws.Range("EA2").Value = -0.03069589261656241
$ws.Range("EB2").Value = 0
$ws.Range("EF2").Value = -0
$ws.Range("EG2").Value = 0.04440782764334671
$ws.Range("EI2").Value = 0.0778287514814467
$ws.Range("EJ2").Value = -0.02568372809262882
$ws.Range("EO2").Value = 0
$ws.Range("EP2").Value = 0.05103484235418234
$ws.Range("EQ2").Value = 0
$ws.Range("ER2").Value = -0.03681046101146412
$ws.Range("ES2").Value = 0.03196449392776059
$ws.Range("ET2").Value = 0
$ws.Range("EU2").Value = -0
$ws.Range("EV2").Value = 0
$ws.Range("EX2").Value = 0
$ws.Range("EY2").Value = 0.04619994784946117
$ws.Range("FA2").Value = -0.03583138226814911
$ws.Range("FB2").Value = 0.0215542075157635
$ws.Range("FD2").Value = -0
$ws.Range("FG2").Value = -0
$ws.Range("FH2").Value = 0.008348946584245088
$ws.Range("FJ2").Value = -0.009127320605065678
$ws.Range("FK2").Value = -0.005088268726051238
$ws.Range("FL2").Value = -0
$ws.Range("FN2").Value = -0
$ws.Range("FP2").Value = -0
$ws.Range("FQ2").Value = -0.004157056812444756
$ws.Range("FR2").Value = -0
$ws.Range("FS2").Value = -0.02284002279882993
$ws.Range("FT2").Value = 0.01101548204690669
$ws.Range("FV2").Value = -0
$ws.Range("FW2").Value = -0
$ws.Range("FY2").Value = 0
$ws.Range("FZ2").Value = -0.03916468553119062
$ws.Range("GB2").Value = 0.02405378880385655
$ws.Range("GD2").Value = 0
$ws.Range("GE2").Value = -0
